# 3CM1_IS_T15_ACME.xlsx - "Plan de riesgos" update
# Adds a new risk row (row 8: "Nuestro cliente no acepte los prototipos"),
# promotes the previously-empty placeholder row into real data, fixes a
# typo in the existing "capacitación" risk text, sets the new row's
# height to match its wrapped two-line content, and moves the active
# selection to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a spelling slip in the existing row 5 risk description ---
# "capacitación" -> "capasitación" (same risk, just a retyped/garbled word)
$ws.Range("D5").Value = "No adquirir la capasitación necesaria para el desarrollo de aplicaciónes móviles "

# --- Populate the new risk entry in row 8 (previously a blank template row) ---
$ws.Range("C8").Value = 42015
$ws.Range("D8").Value = "Nuestro cliente no acepte los prototipos"
$ws.Range("E8").Value = "Desfase de tiempo para rediseñarlos."
$ws.Range("F8").Value = "Media"
$ws.Range("G8").Value = "Baja"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = "Hacer que los prototipos cumplan con la mayoría de los requisitos propuestos por el cliente."
$ws.Range("J8").Value = "Volver a hacer nuevos prototipos."
$ws.Range("K8").Value = "A, T, M, J"

# Row now wraps onto two lines like rows 4/5/7 above it.
$ws.Rows.Item(8).RowHeight = 25.5

# --- Move the active selection, as left by the editor, to E8 ---
[void]$ws.Range("E8").Select()
